$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 86
$ws.Range("F4").Value = 26
$ws.Range("H4").Value = 26
$ws.Range("E5").Value = 119
$ws.Range("F5").Value = 74
$ws.Range("H5").Value = 74
$ws.Range("E10").Value = 398
$ws.Range("F10").Value = 189
$ws.Range("H10").Value = 189
$ws.Range("E11").Value = 272
$ws.Range("F11").Value = 141
$ws.Range("H11").Value = 141
$ws.Range("E12").Value = 386
$ws.Range("E13").Value = 104
$ws.Range("F13").Value = 55
$ws.Range("H13").Value = 55
$ws.Range("E14").Value = 99
$ws.Range("E15").Value = 126
$ws.Range("E17").Value = 75
$ws.Range("E21").Value = 124
$ws.Range("F21").Value = 66
$ws.Range("H21").Value = 66
$ws.Range("E22").Value = 145
$ws.Range("E24").Value = 171
$ws.Range("F24").Value = 87
$ws.Range("H24").Value = 87
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 91
$ws.Range("H25").Value = 91
$ws.Range("E26").Value = 115
$ws.Range("E27").Value = 263
$ws.Range("F27").Value = 125
$ws.Range("H27").Value = 125
$ws.Range("E28").Value = 153
$ws.Range("F28").Value = 51
$ws.Range("H28").Value = 51
$ws.Range("E30").Value = 171
$ws.Range("F30").Value = 95
$ws.Range("H30").Value = 95
$ws.Range("E32").Value = 159
$ws.Range("F32").Value = 90
$ws.Range("H32").Value = 90
$ws.Range("E33").Value = 242
$ws.Range("E35").Value = 115
$ws.Range("E36").Value = 51
$ws.Range("E40").Value = 217
$ws.Range("F40").Value = 97
$ws.Range("H40").Value = 97
$ws.Range("E41").Value = 322
$ws.Range("F41").Value = 143
$ws.Range("H41").Value = 143
$ws.Range("E42").Value = 289
$ws.Range("F42").Value = 152
$ws.Range("H42").Value = 152
$ws.Range("E43").Value = 97
$ws.Range("E44").Value = 255
$ws.Range("F44").Value = 120
$ws.Range("H44").Value = 120
$ws.Range("E45").Value = 113
$ws.Range("F45").Value = 53
$ws.Range("H45").Value = 53
$ws.Range("E46").Value = 252
$ws.Range("F46").Value = 137
$ws.Range("H46").Value = 137
$ws.Range("E47").Value = 362
$ws.Range("F47").Value = 175
$ws.Range("H47").Value = 175
$ws.Range("E48").Value = 170
$ws.Range("F48").Value = 67
$ws.Range("H48").Value = 67
$ws.Range("E49").Value = 243
$ws.Range("E50").Value = 208
$ws.Range("E51").Value = 194
